$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures (columns D and E) for rows 2-51.
# Rows 37 and 38 additionally swap the Algorand / VeChain entries
# (name in column B, link in column C) along with their price/volume data.
#
# Numeric-looking price strings (e.g. "1.006", "0.00001072") are written
# with a leading apostrophe so Excel stores them as text rather than
# auto-converting them into numbers (matching the original inline-string
# formatting, which used values like "20.60" and "0.6450" that would lose
# their trailing zeros / precision if parsed as numbers).

$ws.Range("D2").Value = '27.550.11'
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("D3").Value = '1.757.47'
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").Value = "'324.96"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = "'0.4459"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.3704"
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").Value = "'45.14"
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = "'0.07694"
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").Value = "'1.115"
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = "'21.64"
$ws.Range("E13").Value = '  -4.49%  '
$ws.Range("D14").Value = "'6.157"
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").Value = "'7.427"
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D16").Value = '1.761.02'
$ws.Range("E16").Value = '  -2.29%  '
$ws.Range("D17").Value = "'90.57"
$ws.Range("E17").Value = '  +12.00%  '
$ws.Range("D18").Value = "'0.00001072"
$ws.Range("E18").Value = '  -2.19%  '
$ws.Range("D19").Value = "'0.06276"
$ws.Range("E19").Value = '  -7.64%  '
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = "'17.43"
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Value = "'6.175"
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").Value = "'0.5314"
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("D24").Value = '27.608.56'
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("D25").Value = "'11.56"
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("D26").Value = "'2.307"
$ws.Range("E26").Value = '  -4.20%  '
$ws.Range("D27").Value = "'20.60"
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = "'153.58"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Value = "'2.294"
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("D30").Value = '1.959.55'
$ws.Range("E30").Value = '  -2.26%  '
$ws.Range("D31").Value = "'127.88"
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("D32").Value = "'1.181"
$ws.Range("E32").Value = '  -6.17%  '
$ws.Range("D33").Value = "'5.725"
$ws.Range("E33").Value = '  -2.20%  '
$ws.Range("D34").Value = "'0.09217"
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = "'3.642"
$ws.Range("E35").Value = '  -9.21%  '
$ws.Range("E36").Value = '  +3.47%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.02314"
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2163"
$ws.Range("E38").Value = '  -5.63%  '
$ws.Range("D39").Value = "'0.06107"
$ws.Range("E39").Value = '  -4.24%  '
$ws.Range("D40").Value = "'0.6450"
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("D41").Value = "'5.049"
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("D42").Value = "'1.175"
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("D43").Value = "'7.972"
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").Value = "'1.395"
$ws.Range("E45").Value = '  -4.34%  '
$ws.Range("D46").Value = "'13.65"
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = "'0.5972"
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").Value = "'3.728"
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("D49").Value = "'126.13"
$ws.Range("E49").Value = '  -2.23%  '
$ws.Range("D50").Value = "'1.999"
$ws.Range("E50").Value = '  -2.11%  '
$ws.Range("D51").Value = "'0.06895"
$ws.Range("E51").Value = '  -2.91%  '
